$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append two new rows (106 and 107) with data, matching existing table layout:
# column A = Date (text), B = Terme (text), C = Numero de page (number), D = Occurences (number)
# Column A must stay as plain text (not get auto-converted to a date serial),
# so force the cell format to Text before writing the value.

$ws.Range("A106:A107").NumberFormat = "@"

$ws.Range("A106").Value = "2025-05-09"
$ws.Range("B106").Value = "ruissellement"
$ws.Range("C106").Value = 43
$ws.Range("D106").Value = 1

$ws.Range("A107").Value = "2025-05-09"
$ws.Range("B107").Value = "ruissellement"
$ws.Range("C107").Value = 47
$ws.Range("D107").Value = 2
